$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data was added for "Espárragos" at the top of this
# block. Insert 2 new rows before row 80, shifting the previously-newest
# entries (old rows 80-84) down to rows 82-86.
$ws.Rows.Item(80).Resize(2).Insert()

# Populate the new row 80 with the latest "Primera" quality entry.
$ws.Cells.Item(80, 1).Value2 = 9
$ws.Cells.Item(80, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(80, 3).Value2 = "Metropolitana"
$ws.Cells.Item(80, 4).Value2 = 44504
$ws.Cells.Item(80, 5).Value2 = 13
$ws.Cells.Item(80, 6).Value2 = 300000000
$ws.Cells.Item(80, 7).Value2 = "Espárragos"
$ws.Cells.Item(80, 8).Value2 = "Sin especificar"
$ws.Cells.Item(80, 9).Value2 = "Primera"
$ws.Cells.Item(80, 10).Value2 = 106
$ws.Cells.Item(80, 11).Value2 = 12000
$ws.Cells.Item(80, 12).Value2 = 12000
$ws.Cells.Item(80, 13).Value2 = 12000
$ws.Cells.Item(80, 14).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(80, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(80, 16).Value2 = 1200
$ws.Cells.Item(80, 17).Value2 = 10
$ws.Cells.Item(80, 18).Value2 = "Hortaliza"

# Populate the new row 81 with the latest "Segunda" quality entry.
$ws.Cells.Item(81, 1).Value2 = 9
$ws.Cells.Item(81, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(81, 3).Value2 = "Metropolitana"
$ws.Cells.Item(81, 4).Value2 = 44504
$ws.Cells.Item(81, 5).Value2 = 13
$ws.Cells.Item(81, 6).Value2 = 300000000
$ws.Cells.Item(81, 7).Value2 = "Espárragos"
$ws.Cells.Item(81, 8).Value2 = "Sin especificar"
$ws.Cells.Item(81, 9).Value2 = "Segunda"
$ws.Cells.Item(81, 10).Value2 = 52
$ws.Cells.Item(81, 11).Value2 = 10000
$ws.Cells.Item(81, 12).Value2 = 10000
$ws.Cells.Item(81, 13).Value2 = 10000
$ws.Cells.Item(81, 14).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(81, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(81, 16).Value2 = 1000
$ws.Cells.Item(81, 17).Value2 = 10
$ws.Cells.Item(81, 18).Value2 = "Hortaliza"
